{"js": "const replacements = [\n  [\"2025-05-26 Monday\", \"2025-05-27 Tuesday\"],\n  [\"979\u00f77=139, 6\", \"999\u00f73=333, 0\"],\n  [\"346\u00f79=38, 4\", \"337\u00f79=37, 4\"],\n  [\"293\u00f76=48, 5\", \"308\u00f72=154, 0\"],\n  [\"219\u00f76=36, 3\", \"535\u00f76=89, 1\"],\n  [\"615\u00f74=153, 3\", \"198\u00f76=33, 0\"],\n  [\"746\u00f74=186, 2\", \"579\u00f79=64, 3\"],\n  [\"664\u00f79=73, 7\", \"724\u00f79=80, 4\"],\n  [\"923\u00f77=131, 6\", \"429\u00f76=71, 3\"],\n  [\"615\u00f73=205, 0\", \"459\u00f77=65, 4\"],\n  [\"276\u00f74=69, 0\", \"766\u00f79=85, 1\"],\n  [\"403\u00f74=100, 3\", \"322\u00f73=107, 1\"],\n  [\"695\u00f73=231, 2\", \"311\u00f73=103, 2\"],\n  [\"637\u00f72=318, 1\", \"388\u00f79=43, 1\"],\n  [\"766\u00f75=153, 1\", \"309\u00f78=38, 5\"],\n  [\"153\u00f73=51, 0\", \"599\u00f78=74, 7\"],\n  [\"549\u00f76=91, 3\", \"982\u00f74=245, 2\"],\n  [\"651\u00f72=325, 1\", \"877\u00f77=125, 2\"],\n  [\"104\u00f72=52, 0\", \"972\u00f74=243, 0\"],\n  [\"496\u00f78=62, 0\", \"622\u00f78=77, 6\"],\n  [\"633\u00f77=90, 3\", \"843\u00f73=281, 0\"],\n  [\"594\u00f78=74, 2\", \"653\u00f79=72, 5\"],\n  [\"376\u00f78=47, 0\", \"965\u00f74=241, 1\"],\n  [\"901\u00f78=112, 5\", \"983\u00f73=327, 2\"],\n  [\"262\u00f78=32, 6\", \"476\u00f73=158, 2\"],\n  [\"264\u00f76=44, 0\", \"628\u00f78=78, 4\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load text for each paragraph's runs by re-loading ranges.\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const text = p.text;\n  for (const [oldText, newText] of replacements) {\n    if (text === oldText) {\n      const range = p.getRange(\"Whole\");\n      range.insertText(newText, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-26 Monday\", \"2025-05-27 Tuesday\"),\n    @(\"979\u00f77=139, 6\", \"999\u00f73=333, 0\"),\n    @(\"346\u00f79=38, 4\", \"337\u00f79=37, 4\"),\n    @(\"293\u00f76=48, 5\", \"308\u00f72=154, 0\"),\n    @(\"219\u00f76=36, 3\", \"535\u00f76=89, 1\"),\n    @(\"615\u00f74=153, 3\", \"198\u00f76=33, 0\"),\n    @(\"746\u00f74=186, 2\", \"579\u00f79=64, 3\"),\n    @(\"664\u00f79=73, 7\", \"724\u00f79=80, 4\"),\n    @(\"923\u00f77=131, 6\", \"429\u00f76=71, 3\"),\n    @(\"615\u00f73=205, 0\", \"459\u00f77=65, 4\"),\n    @(\"276\u00f74=69, 0\", \"766\u00f79=85, 1\"),\n    @(\"403\u00f74=100, 3\", \"322\u00f73=107, 1\"),\n    @(\"695\u00f73=231, 2\", \"311\u00f73=103, 2\"),\n    @(\"637\u00f72=318, 1\", \"388\u00f79=43, 1\"),\n    @(\"766\u00f75=153, 1\", \"309\u00f78=38, 5\"),\n    @(\"153\u00f73=51, 0\", \"599\u00f78=74, 7\"),\n    @(\"549\u00f76=91, 3\", \"982\u00f74=245, 2\"),\n    @(\"651\u00f72=325, 1\", \"877\u00f77=125, 2\"),\n    @(\"104\u00f72=52, 0\", \"972\u00f74=243, 0\"),\n    @(\"496\u00f78=62, 0\", \"622\u00f78=77, 6\"),\n    @(\"633\u00f77=90, 3\", \"843\u00f73=281, 0\"),\n    @(\"594\u00f78=74, 2\", \"653\u00f79=72, 5\"),\n    @(\"376\u00f78=47, 0\", \"965\u00f74=241, 1\"),\n    @(\"901\u00f78=112, 5\", \"983\u00f73=327, 2\"),\n    @(\"262\u00f78=32, 6\", \"476\u00f73=158, 2\"),\n    @(\"264\u00f76=44, 0\", \"628\u00f78=78, 4\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    [void]$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
